$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "64.335.11"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "3.486.92"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  -0.18%  "
Set-TextValue $ws "D5" "587.03"
$ws.Range("E5").Value = "  +0.73%  "
Set-TextValue $ws "D6" "134.30"
$ws.Range("E6").Value = "  +2.43%  "
$ws.Range("D7").Value = "3.487.54"
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("E10").Value = "  -0.03%  "
Set-TextValue $ws "D11" "7.20"
$ws.Range("E11").Value = "  -0.13%  "
Set-TextValue $ws "D12" "0.376"
$ws.Range("E12").Value = "  -2.50%  "
$ws.Range("D13").Value = "4.079.50"
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("E14").Value = "  +1.69%  "
Set-TextValue $ws "D15" "0.0000180"
$ws.Range("E15").Value = "  +1.34%  "
$ws.Range("D16").Value = "3.483.32"
$ws.Range("D17").Value = "64.362.24"
$ws.Range("E17").Value = "  +0.21%  "
Set-TextValue $ws "D18" "25.13"
$ws.Range("E18").Value = "  -9.18%  "
Set-TextValue $ws "D19" "9.95"
$ws.Range("E19").Value = "  +0.29%  "
Set-TextValue $ws "D20" "5.66"
$ws.Range("E20").Value = "  +0.22%  "
Set-TextValue $ws "D21" "13.74"
$ws.Range("E21").Value = "  -3.44%  "
Set-TextValue $ws "D22" "385.84"
$ws.Range("E22").Value = "  -1.68%  "
$ws.Range("E23").Value = "  -1.57%  "
$ws.Range("D24").Value = "3.624.20"
$ws.Range("E24").Value = "  +0.11%  "
Set-TextValue $ws "D25" "74.02"
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  +5.37%  "
$ws.Range("E29").Value = "  -0.16%  "
Set-TextValue $ws "D30" "0.999"
$ws.Range("E30").Value = "  +0.39%  "
Set-TextValue $ws "D31" "7.43"
$ws.Range("E31").Value = "  -0.13%  "
Set-TextValue $ws "D32" "2.23"
$ws.Range("E32").Value = "  -0.69%  "
Set-TextValue $ws "D33" "8.20"
$ws.Range("E33").Value = "  +0.48%  "
$ws.Range("D34").Value = "3.506.63"
$ws.Range("E34").Value = "  +0.79%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  +1.66%  "
Set-TextValue $ws "D37" "23.39"
$ws.Range("E37").Value = "  -1.85%  "
Set-TextValue $ws "D38" "5.29"
$ws.Range("E38").Value = "  +1.01%  "
Set-TextValue $ws "D39" "6.83"
$ws.Range("E39").Value = "  -1.83%  "
$ws.Range("E40").Value = "  -2.17%  "
Set-TextValue $ws "D41" "162.41"
$ws.Range("E41").Value = "  -4.45%  "
Set-TextValue $ws "D42" "0.0778"
$ws.Range("E42").Value = "  -3.06%  "
Set-TextValue $ws "D43" "0.803"
$ws.Range("E43").Value = "  -0.93%  "
Set-TextValue $ws "D44" "25.55"
$ws.Range("E44").Value = "  +0.54%  "
$ws.Range("E45").Value = "  -0.12%  "
Set-TextValue $ws "D46" "41.75"
$ws.Range("E46").Value = "  +0.21%  "
Set-TextValue $ws "D47" "4.39"
$ws.Range("E47").Value = "  +0.99%  "
$ws.Range("E48").Value = "  +0.62%  "
$ws.Range("E49").Value = "  +2.01%  "
$ws.Range("D50").Value = "2.470.04"
$ws.Range("E50").Value = "  +2.09%  "
$ws.Range("E51").Value = "  -1.66%  "